# Generate Report for Handback
#
# The localization round-trip for both target files (7f90948e... and
# c38e6498...) has completed for zh-cn and de-de. Update the per-language
# status sheets with:
#   - Status -> "Handed back: in sync with en-US"
#   - Latest Target File  (col I) -> source .md file, hyperlinked
#   - Latest Handback File (col J) -> the generated xlf for that language
#   - Latest Handback DateTime (col K) -> timestamp of the handback
#
# (The Overview sheet's Status/columns derive from the same shared
#  strings, so it updates automatically.)

$wb = $excel.ActiveWorkbook

$fileA = "7f90948e-49c0-4729-b772-4af0c8b172e1.md"
$fileB = "c38e6498-941e-40b3-963b-9971120a447b.md"
$urlA  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71ff91d65a0d3c883ef6be43eb99950d41b8ad81/e2e/$fileA"
$urlB  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71ff91d65a0d3c883ef6be43eb99950d41b8ad81/e2e/$fileB"

$statusText = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = "7f90948e-49c0-4729-b772-4af0c8b172e1.0497c110ada34583356bed85bd6ebe3831562394.zh-cn.xlf"
$wsZh.Range("J3").Value = "c38e6498-941e-40b3-963b-9971120a447b.9f0fdb17d654e5a9d4962fcd2d715f6351e7b96d.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-29 20:42:08"
$wsZh.Range("K3").Value = "2016-08-29 20:42:08"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", $fileA)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlB, "", "", $fileB)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = "7f90948e-49c0-4729-b772-4af0c8b172e1.0497c110ada34583356bed85bd6ebe3831562394.de-de.xlf"
$wsDe.Range("J3").Value = "c38e6498-941e-40b3-963b-9971120a447b.9f0fdb17d654e5a9d4962fcd2d715f6351e7b96d.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-29 20:42:18"
$wsDe.Range("K3").Value = "2016-08-29 20:42:18"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", $fileA)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlB, "", "", $fileB)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

# --- Overview sheet: status + columns get wider too ---------------------
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("E2").Value = $statusText
$wsOv.Range("F2").Value = $statusText
$wsOv.Range("E3").Value = $statusText
$wsOv.Range("F3").Value = $statusText
$wsOv.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOv.Columns.Item(6).ColumnWidth = 29.9777047293527

Write-Host "Handback report generated."
